# Update TLKGY yearly financials with refreshed figures (income statement,
# balance sheet and cash flow sections) pulled from the latest data run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 4).Value = 2811400
$ws.Cells.Item(8, 5).Value = 2808100
$ws.Cells.Item(8, 6).Value = 2558300
$ws.Cells.Item(8, 7).Value = 2245400
$ws.Cells.Item(8, 8).Value = 2144500
$ws.Cells.Item(8, 9).Value = 2203100
$ws.Cells.Item(8, 10).Value = 2267200
$ws.Cells.Item(9, 4).Value = 607400
$ws.Cells.Item(9, 5).Value = 624800
$ws.Cells.Item(9, 6).Value = 534900
$ws.Cells.Item(9, 7).Value = 423500
$ws.Cells.Item(9, 8).Value = 419500
$ws.Cells.Item(9, 9).Value = 411900
$ws.Cells.Item(9, 10).Value = 375900
$ws.Cells.Item(10, 4).Value = 2204000
$ws.Cells.Item(10, 5).Value = 2183300
$ws.Cells.Item(10, 6).Value = 2023400
$ws.Cells.Item(10, 7).Value = 1821900
$ws.Cells.Item(10, 8).Value = 1724900
$ws.Cells.Item(10, 9).Value = 1791200
$ws.Cells.Item(10, 10).Value = 1891400
$ws.Cells.Item(14, 4).Value = 3200
$ws.Cells.Item(14, 5).Value = 9800
$ws.Cells.Item(14, 6).Value = 11700
$ws.Cells.Item(14, 7).Value = 15100
$ws.Cells.Item(14, 8).Value = 43300
$ws.Cells.Item(14, 9).Value = 843700
$ws.Cells.Item(14, 10).Value = 55900
$ws.Cells.Item(15, 4).Value = 380900
$ws.Cells.Item(15, 5).Value = 378200
$ws.Cells.Item(15, 6).Value = 365200
$ws.Cells.Item(15, 7).Value = 362200
$ws.Cells.Item(15, 8).Value = 360500
$ws.Cells.Item(15, 9).Value = 759100
$ws.Cells.Item(15, 10).Value = 364300
$ws.Cells.Item(17, 4).Value = 2472900
$ws.Cells.Item(17, 5).Value = 2450700
$ws.Cells.Item(17, 6).Value = 2333600
$ws.Cells.Item(17, 7).Value = 2016700
$ws.Cells.Item(17, 8).Value = 1865100
$ws.Cells.Item(17, 9).Value = 2967000
$ws.Cells.Item(17, 10).Value = 2102200
$ws.Cells.Item(18, 4).Value = 338500
$ws.Cells.Item(18, 5).Value = 357400
$ws.Cells.Item(18, 6).Value = 224700
$ws.Cells.Item(18, 7).Value = 228600
$ws.Cells.Item(18, 8).Value = 279400
$ws.Cells.Item(18, 9).Value = -763900
$ws.Cells.Item(18, 10).Value = 165000
$ws.Cells.Item(20, 4).Value = 16800
$ws.Cells.Item(20, 5).Value = -3500
$ws.Cells.Item(20, 6).Value = 7000
$ws.Cells.Item(20, 7).Value = -12300
$ws.Cells.Item(20, 8).Value = -8800
$ws.Cells.Item(20, 9).Value = 49000
$ws.Cells.Item(20, 10).Value = -59600
$ws.Cells.Item(21, 4).Value = 739500
$ws.Cells.Item(21, 5).Value = 741900
$ws.Cells.Item(21, 6).Value = 610500
$ws.Cells.Item(21, 7).Value = 593600
$ws.Cells.Item(21, 8).Value = 677400
$ws.Cells.Item(21, 9).Value = 538900
$ws.Cells.Item(21, 10).Value = 526200
$ws.Cells.Item(22, 4).Value = 61200
$ws.Cells.Item(22, 5).Value = 42400
$ws.Cells.Item(22, 6).Value = 35700
$ws.Cells.Item(22, 9).Value = 45200
$ws.Cells.Item(22, 10).Value = 52400
$ws.Cells.Item(23, 4).Value = 294100
$ws.Cells.Item(23, 5).Value = 311500
$ws.Cells.Item(23, 6).Value = 196000
$ws.Cells.Item(23, 7).Value = 216300
$ws.Cells.Item(23, 8).Value = 270500
$ws.Cells.Item(23, 9).Value = -760100
$ws.Cells.Item(23, 10).Value = 53000
$ws.Cells.Item(24, 4).Value = 77700
$ws.Cells.Item(24, 5).Value = 47400
$ws.Cells.Item(24, 6).Value = 36900
$ws.Cells.Item(24, 7).Value = -1900
$ws.Cells.Item(24, 8).Value = 24500
$ws.Cells.Item(24, 9).Value = 30300
$ws.Cells.Item(24, 10).Value = 40800
$ws.Cells.Item(26, 4).Value = 216400
$ws.Cells.Item(26, 5).Value = 264200
$ws.Cells.Item(26, 6).Value = 159100
$ws.Cells.Item(26, 7).Value = 218200
$ws.Cells.Item(26, 8).Value = 246100
$ws.Cells.Item(26, 9).Value = -790400
$ws.Cells.Item(26, 10).Value = 12300
$ws.Cells.Item(27, 4).Value = 209200
$ws.Cells.Item(27, 5).Value = 260200
$ws.Cells.Item(27, 6).Value = 151500
$ws.Cells.Item(27, 7).Value = 211000
$ws.Cells.Item(27, 8).Value = 237800
$ws.Cells.Item(27, 9).Value = -798800
$ws.Cells.Item(27, 10).Value = 3600
$ws.Cells.Item(29, 7).Value = 25200
$ws.Cells.Item(29, 8).Value = 24200
$ws.Cells.Item(29, 9).Value = -7100
$ws.Cells.Item(29, 10).Value = -18400
$ws.Cells.Item(32, 4).Value = -16800
$ws.Cells.Item(32, 5).Value = 3500
$ws.Cells.Item(32, 6).Value = -7000
$ws.Cells.Item(32, 7).Value = 12300
$ws.Cells.Item(32, 8).Value = 8800
$ws.Cells.Item(32, 9).Value = -49000
$ws.Cells.Item(32, 10).Value = 59600
$ws.Cells.Item(33, 4).Value = 209200
$ws.Cells.Item(33, 5).Value = 260200
$ws.Cells.Item(33, 6).Value = 151500
$ws.Cells.Item(33, 7).Value = 236200
$ws.Cells.Item(33, 8).Value = 262000
$ws.Cells.Item(33, 9).Value = -806000
$ws.Cells.Item(33, 10).Value = -14800
$ws.Cells.Item(35, 4).Value = 209200
$ws.Cells.Item(35, 5).Value = 260200
$ws.Cells.Item(35, 6).Value = 151500
$ws.Cells.Item(35, 7).Value = 236200
$ws.Cells.Item(35, 8).Value = 262000
$ws.Cells.Item(35, 9).Value = -806000
$ws.Cells.Item(35, 10).Value = -14800
$ws.Cells.Item(41, 4).Value = 114400
$ws.Cells.Item(41, 5).Value = 65300
$ws.Cells.Item(41, 6).Value = 28600
$ws.Cells.Item(41, 7).Value = 11100
$ws.Cells.Item(41, 8).Value = 13200
$ws.Cells.Item(41, 9).Value = 15800
$ws.Cells.Item(41, 10).Value = 48900
$ws.Cells.Item(42, 4).Value = 176000
$ws.Cells.Item(42, 5).Value = 208800
$ws.Cells.Item(42, 6).Value = 146000
$ws.Cells.Item(42, 7).Value = 238600
$ws.Cells.Item(42, 8).Value = 113000
$ws.Cells.Item(42, 9).Value = 575800
$ws.Cells.Item(42, 10).Value = 181600
$ws.Cells.Item(43, 4).Value = 1097300
$ws.Cells.Item(43, 5).Value = 534800
$ws.Cells.Item(43, 6).Value = 410200
$ws.Cells.Item(43, 7).Value = 311800
$ws.Cells.Item(43, 8).Value = 302100
$ws.Cells.Item(43, 9).Value = 1029300
$ws.Cells.Item(43, 10).Value = 328900
$ws.Cells.Item(44, 4).Value = 196700
$ws.Cells.Item(44, 5).Value = 94900
$ws.Cells.Item(44, 6).Value = 133100
$ws.Cells.Item(44, 7).Value = 43700
$ws.Cells.Item(44, 8).Value = 44300
$ws.Cells.Item(44, 9).Value = 104200
$ws.Cells.Item(44, 10).Value = 68100
$ws.Cells.Item(45, 4).Value = 36100
$ws.Cells.Item(45, 5).Value = 49700
$ws.Cells.Item(45, 6).Value = 230300
$ws.Cells.Item(45, 7).Value = 157400
$ws.Cells.Item(45, 8).Value = 100800
$ws.Cells.Item(45, 9).Value = 184200
$ws.Cells.Item(45, 10).Value = 72100
$ws.Cells.Item(46, 4).Value = 965200
$ws.Cells.Item(46, 5).Value = 953500
$ws.Cells.Item(46, 6).Value = 881700
$ws.Cells.Item(46, 7).Value = 762600
$ws.Cells.Item(46, 8).Value = 573400
$ws.Cells.Item(46, 9).Value = 769200
$ws.Cells.Item(46, 10).Value = 699500
$ws.Cells.Item(47, 4).Value = 28900
$ws.Cells.Item(47, 5).Value = 28100
$ws.Cells.Item(47, 6).Value = 181900
$ws.Cells.Item(47, 7).Value = 185100
$ws.Cells.Item(47, 8).Value = 208000
$ws.Cells.Item(47, 9).Value = 197700
$ws.Cells.Item(47, 10).Value = 174900
$ws.Cells.Item(48, 4).Value = 4164100
$ws.Cells.Item(48, 5).Value = 1913500
$ws.Cells.Item(48, 6).Value = 1925900
$ws.Cells.Item(48, 7).Value = 1677800
$ws.Cells.Item(48, 8).Value = 1721900
$ws.Cells.Item(48, 10).Value = 2478100
$ws.Cells.Item(49, 4).Value = 615800
$ws.Cells.Item(49, 5).Value = 323500
$ws.Cells.Item(49, 6).Value = 301900
$ws.Cells.Item(49, 7).Value = 204400
$ws.Cells.Item(49, 8).Value = 194200
$ws.Cells.Item(49, 9).Value = 530700
$ws.Cells.Item(49, 10).Value = 243700
$ws.Cells.Item(52, 4).Value = 77200
$ws.Cells.Item(52, 5).Value = 74600
$ws.Cells.Item(52, 6).Value = 87700
$ws.Cells.Item(52, 7).Value = 112300
$ws.Cells.Item(52, 8).Value = 3300
$ws.Cells.Item(52, 9).Value = 6200
$ws.Cells.Item(52, 10).Value = 6900
$ws.Cells.Item(54, 4).Value = 3461300
$ws.Cells.Item(54, 5).Value = 3293300
$ws.Cells.Item(54, 6).Value = 3190700
$ws.Cells.Item(54, 7).Value = 2877400
$ws.Cells.Item(54, 8).Value = 2700800
$ws.Cells.Item(54, 9).Value = 2849600
$ws.Cells.Item(54, 10).Value = 3603000
$ws.Cells.Item(57, 4).Value = 224300
$ws.Cells.Item(57, 5).Value = 265200
$ws.Cells.Item(57, 6).Value = 489000
$ws.Cells.Item(57, 7).Value = 386200
$ws.Cells.Item(57, 8).Value = 350900
$ws.Cells.Item(57, 9).Value = 319300
$ws.Cells.Item(57, 10).Value = 194600
$ws.Cells.Item(58, 4).Value = 156100
$ws.Cells.Item(58, 5).Value = 112000
$ws.Cells.Item(58, 6).Value = 48600
$ws.Cells.Item(58, 7).Value = 110600
$ws.Cells.Item(58, 8).Value = 22100
$ws.Cells.Item(58, 9).Value = 378300
$ws.Cells.Item(58, 10).Value = 88600
$ws.Cells.Item(59, 4).Value = 505600
$ws.Cells.Item(59, 5).Value = 523300
$ws.Cells.Item(59, 6).Value = 359200
$ws.Cells.Item(59, 7).Value = 302900
$ws.Cells.Item(59, 8).Value = 319400
$ws.Cells.Item(59, 9).Value = 522000
$ws.Cells.Item(59, 10).Value = 382300
$ws.Cells.Item(60, 4).Value = 886000
$ws.Cells.Item(60, 5).Value = 900500
$ws.Cells.Item(60, 6).Value = 896800
$ws.Cells.Item(60, 7).Value = 787000
$ws.Cells.Item(60, 8).Value = 692300
$ws.Cells.Item(60, 9).Value = 899800
$ws.Cells.Item(60, 10).Value = 665500
$ws.Cells.Item(61, 4).Value = 491100
$ws.Cells.Item(61, 5).Value = 325200
$ws.Cells.Item(61, 6).Value = 313000
$ws.Cells.Item(61, 7).Value = 222300
$ws.Cells.Item(61, 8).Value = 258700
$ws.Cells.Item(61, 9).Value = 267200
$ws.Cells.Item(61, 10).Value = 404200
$ws.Cells.Item(62, 4).Value = 213400
$ws.Cells.Item(62, 5).Value = 154900
$ws.Cells.Item(62, 6).Value = 174000
$ws.Cells.Item(62, 7).Value = 139000
$ws.Cells.Item(62, 8).Value = 163200
$ws.Cells.Item(62, 9).Value = 790500
$ws.Cells.Item(62, 10).Value = 467500
$ws.Cells.Item(66, 4).Value = 1615100
$ws.Cells.Item(66, 5).Value = 1403700
$ws.Cells.Item(66, 6).Value = 1410400
$ws.Cells.Item(66, 7).Value = 1173300
$ws.Cells.Item(66, 8).Value = 1140100
$ws.Cells.Item(66, 9).Value = 1629700
$ws.Cells.Item(66, 10).Value = 1566900
$ws.Cells.Item(72, 4).Value = 1500100
$ws.Cells.Item(72, 5).Value = 1532600
$ws.Cells.Item(72, 6).Value = 1423400
$ws.Cells.Item(72, 7).Value = 1347200
$ws.Cells.Item(72, 8).Value = 1256600
$ws.Cells.Item(72, 9).Value = 1064100
$ws.Cells.Item(72, 10).Value = 1732000
$ws.Cells.Item(76, 4).Value = 1846200
$ws.Cells.Item(76, 5).Value = 1889600
$ws.Cells.Item(76, 6).Value = 1780300
$ws.Cells.Item(76, 7).Value = 1704200
$ws.Cells.Item(76, 8).Value = 1560700
$ws.Cells.Item(76, 9).Value = 1219900
$ws.Cells.Item(76, 10).Value = 2036100
$ws.Cells.Item(81, 4).Value = 209200
$ws.Cells.Item(81, 5).Value = 260200
$ws.Cells.Item(81, 6).Value = 151500
$ws.Cells.Item(81, 7).Value = 236200
$ws.Cells.Item(81, 8).Value = 262000
$ws.Cells.Item(81, 9).Value = -806000
$ws.Cells.Item(81, 10).Value = -14800
$ws.Cells.Item(83, 4).Value = 384200
$ws.Cells.Item(83, 5).Value = 388000
$ws.Cells.Item(83, 6).Value = 378800
$ws.Cells.Item(83, 7).Value = 377300
$ws.Cells.Item(83, 8).Value = 406900
$ws.Cells.Item(83, 9).Value = 1253800
$ws.Cells.Item(83, 10).Value = 420700
$ws.Cells.Item(89, 4).Value = 417000
$ws.Cells.Item(89, 5).Value = 379800
$ws.Cells.Item(89, 6).Value = 459300
$ws.Cells.Item(89, 7).Value = 430500
$ws.Cells.Item(89, 8).Value = 436300
$ws.Cells.Item(89, 9).Value = 512100
$ws.Cells.Item(89, 10).Value = 403800
$ws.Cells.Item(91, 4).Value = -532800
$ws.Cells.Item(91, 5).Value = -581200
$ws.Cells.Item(91, 6).Value = -403800
$ws.Cells.Item(91, 7).Value = -347500
$ws.Cells.Item(91, 8).Value = -436600
$ws.Cells.Item(91, 9).Value = -385700
$ws.Cells.Item(91, 10).Value = -320400
$ws.Cells.Item(94, 4).Value = -454700
$ws.Cells.Item(94, 5).Value = -454900
$ws.Cells.Item(94, 6).Value = -563100
$ws.Cells.Item(94, 7).Value = -354200
$ws.Cells.Item(94, 8).Value = -297000
$ws.Cells.Item(94, 9).Value = -378300
$ws.Cells.Item(94, 10).Value = -336300
$ws.Cells.Item(100, 4).Value = 118500
$ws.Cells.Item(100, 5).Value = 4700
$ws.Cells.Item(100, 6).Value = 28200
$ws.Cells.Item(100, 7).Value = 46900
$ws.Cells.Item(100, 8).Value = -177000
$ws.Cells.Item(100, 9).Value = -50100
$ws.Cells.Item(100, 10).Value = -108700
$ws.Cells.Item(101, 9).Value = -300
$ws.Cells.Item(102, 4).Value = 80800
$ws.Cells.Item(102, 5).Value = -70100
$ws.Cells.Item(102, 6).Value = -75400
$ws.Cells.Item(102, 7).Value = 123400
$ws.Cells.Item(102, 8).Value = -37000
$ws.Cells.Item(102, 9).Value = 83400
$ws.Cells.Item(102, 10).Value = -41700
